$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 141804.95
$ws.Range("E2").Value = 162843.23
$ws.Range("F2").Value = 331318.89
$ws.Range("G2").Value = 403018.76

$ws.Range("D3").Value = 150166.71
$ws.Range("E3").Value = 187284.72
$ws.Range("F3").Value = 343648.67
$ws.Range("G3").Value = 456271.58

$ws.Range("D4").Value = 166620.23
$ws.Range("E4").Value = 193786.84
$ws.Range("F4").Value = 323995.99
$ws.Range("G4").Value = 381373.46

$ws.Range("C5").Value = 46161.32
$ws.Range("D5").Value = 130348.58
$ws.Range("E5").Value = 215439.97
$ws.Range("F5").Value = 396408.69
$ws.Range("G5").Value = 424615.57

$ws.Range("D6").Value = 115923.13
$ws.Range("E6").Value = 220807.59
$ws.Range("F6").Value = 378480.22
$ws.Range("G6").Value = 490874.2

$ws.Range("D7").Value = 141324.93
$ws.Range("E7").Value = 246540.72
$ws.Range("F7").Value = 403893.3
$ws.Range("G7").Value = 116449.89

$ws.Range("D8").Value = 136278.09
$ws.Range("E8").Value = 225495.54
$ws.Range("F8").Value = 389554.43

$ws.Range("C9").Value = 117644.99
$ws.Range("D9").Value = 148973.56
$ws.Range("E9").Value = 267326.16
$ws.Range("F9").Value = 382184.28

$ws.Range("D10").Value = 153315.28
$ws.Range("E10").Value = 292751.94
$ws.Range("F10").Value = 387682.65

$ws.Range("C11").Value = 120207.54
$ws.Range("D11").Value = 134867.6
$ws.Range("E11").Value = 332005.29
$ws.Range("F11").Value = 418223.63

$ws.Range("C12").Value = 129704.99
$ws.Range("D12").Value = 151314.9
$ws.Range("E12").Value = 249661.6
$ws.Range("F12").Value = 309720.57

$ws.Range("C13").Value = 106870.37
$ws.Range("E13").Value = 251259.08
$ws.Range("F13").Value = 382246.32
